$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

$ws.Range("B2").Value = 0.701460474100761
$ws.Range("C2").Value = 0.7014221224979839
$ws.Range("D2").Value = 0.00921802781522274
$ws.Range("E2").Value = 0.002188196638599038
$ws.Range("B3").Value = 0.7608322267528042
$ws.Range("C3").Value = 0.8101514534079264
$ws.Range("D3").Value = 0.01531793642789125
$ws.Range("E3").Value = 0.01246019825339317
$ws.Range("B4").Value = 0.6583490089977085
$ws.Range("C4").Value = 0.6551823223469332
$ws.Range("D4").Value = 0.01853840053081512
$ws.Range("E4").Value = 0.003274488728493452
$ws.Range("B5").Value = 0.775545979630454
$ws.Range("C5").Value = 0.7963016222642612
$ws.Range("D5").Value = 0.01077786833047867
$ws.Range("E5").Value = 0.006051451433449984
$ws.Range("B6").Value = 0.63322110625362
$ws.Range("C6").Value = 0.5960019331165983
$ws.Range("D6").Value = 0.01873300410807133
$ws.Range("E6").Value = 0.01038294844329357
$ws.Range("B7").Value = 0.696835444845756
$ws.Range("C7").Value = 0.6772249883860954
$ws.Range("D7").Value = 0.01883435994386673
$ws.Range("E7").Value = 0.0156498197466135
$ws.Range("B8").Value = 0.7031698175963234
$ws.Range("C8").Value = 0.6480114350028756
$ws.Range("D8").Value = 0.01851365901529789
$ws.Range("E8").Value = 0.009519393555819988
$ws.Range("B9").Value = 0.5074973337885037
$ws.Range("C9").Value = 0.5618052683752099
$ws.Range("D9").Value = 0.02887573838233948
$ws.Range("E9").Value = 0.02052079513669014
$ws.Range("B10").Value = 0.7591028727087614
$ws.Range("C10").Value = 0.7269015622091259
$ws.Range("D10").Value = 0.0450722724199295
$ws.Range("E10").Value = 0.04082705825567245
$ws.Range("B11").Value = 0.7894464926494161
$ws.Range("C11").Value = 0.7435859500757549
$ws.Range("D11").Value = 0.02352042868733406
$ws.Range("E11").Value = 0.0149396974593401
$ws.Range("B12").Value = 0.4739061584753052
$ws.Range("C12").Value = 0.4480116858024746
$ws.Range("D12").Value = 0.007459850050508976
$ws.Range("E12").Value = 0.005506517831236124
$ws.Range("B13").Value = 0.6935747003736874
$ws.Range("C13").Value = 0.6640252545140365
$ws.Range("D13").Value = 0.02653098106384277
$ws.Range("E13").Value = 0.01588558219373226
$ws.Range("B14").Value = 0.6997655701428468
$ws.Range("C14").Value = 0.7547041824631651
$ws.Range("D14").Value = 0.02889014780521393
$ws.Range("E14").Value = 0.02135049551725388
$ws.Range("B15").Value = 0.7652468767042275
$ws.Range("C15").Value = 0.7112876378675204
$ws.Range("D15").Value = 0.03103909641504288
$ws.Range("E15").Value = 0.01275253854691982
$ws.Range("B16").Value = 0.8639562914236636
$ws.Range("C16").Value = 0.838347845800362
$ws.Range("D16").Value = 0.02276878617703915
$ws.Range("E16").Value = 0.01191858761012554
$ws.Range("B17").Value = 0.6003726431216587
$ws.Range("C17").Value = 0.6484345613257172
$ws.Range("D17").Value = 0.01317205093801022
$ws.Range("E17").Value = 0.01011515501886606
$ws.Range("B18").Value = 0.7469573475915958
$ws.Range("C18").Value = 0.6616749005906198
$ws.Range("D18").Value = 0.1047332286834717
$ws.Range("E18").Value = 0.07370149344205856
$ws.Range("B19").Value = 0.7185787508686599
$ws.Range("C19").Value = 0.6465363808979195
$ws.Range("D19").Value = 0.03013242594897747
$ws.Range("E19").Value = 0.02225074544548988
